$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column C (STATUS), shifting it to D.
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "NOME_DO_CURSO"

# Re-apply the (Normal) cell style across the header row so the cells pick up
# an explicit style record, matching the updated style definitions.
$ws.Range("A1:D1").Style = "Normal"

# Give the new column a width similar to its neighbours.
$ws.Columns("C:C").ColumnWidth = 16.98

# Move the active selection, as recorded in the saved view state.
$ws.Range("C5").Select() | Out-Null
